$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "Unknown" -> "unknown" for cells D2:J2
$ws.Range("D2:J2").Value = "unknown"
